$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a D column with the character-length of each word in column A,
# for every data row (2 through 82).
$lastRow = 82
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("D$r").Formula = "=LEN(A$r)"
}

# Update the AutoFilter range to include the new column.
$ws.Range("A1:D$lastRow").AutoFilter() | Out-Null

# Re-sort the data (now including the header) first by the new D column,
# then by the word column A.
$sortRange = $ws.Range("A1:D$lastRow")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D2:D$lastRow")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("A2:A$lastRow")) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Update the current selection.
$ws.Range("S8").Select() | Out-Null

# Set up the page for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
